$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(45940, 5385.51843233239, 4708.94040472918, 6956, 8289.083544, 27.3543965165326),
    @(45941, 1407.58714431483, 3121.89434285284, 3620, 3295.322539, 57.9012390640835),
    @(45942, 1273.73040983078, 3143.03816550115, 3620, 3116.796408, 56.9210068195988),
    @(45943, 5633.42878264471, 5507.60551254741, 3620, 8222.527712, 186.529351745946),
    @(45944, 5584.41744308704, 5749.75180121496, 3620, 8131.752676, 194.878626421996),
    @(45945, 5841.28374289435, 6016.21774721853, 3620, 8971.011103, 230.247712805174),
    @(45946, 5841.28374289435, 6277.67514397036, 3620, 8971.011103, 241.141771003167),
    @(45947, 5841.28374289435, 5468.1137657723, 3620, 8971.011103, 207.410046911581),
    @(45948, 1742.27770790123, 3540.60017441242, 3620, 4473.856108, 110.507440604633),
    @(45949, 1636.94065696827, 3610.38316744489, 3620, 4358.695007, 113.005729894859),
    @(45950, 6392.95297294923, 6096.25648379392, 3620, 9928.52543, 250.492872535195),
    @(45951, 6392.95297294923, 6301.52024840312, 3620, 9928.52543, 259.045529393912),
    @(45952, 6392.95297294923, 6175.72572942504, 3620, 9928.52543, 253.804091103158),
    @(45953, 6392.95297294923, 6181.177807063, 3620, 9928.52543, 254.03126100474)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}
